# Update the random data grid (A1:J20) on Sheet1 with the new set of
# generated values, per the commit: "Creato il calendario con la
# separazione delle settimane da gennaio a fine aprile" edit.
# Conditional-formatting color scales on the sheet are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 174.614824342828
$ws.Cells.Item(1, 2).Value = 79.9171846732112
$ws.Cells.Item(1, 3).Value = 163.248516695224
$ws.Cells.Item(1, 4).Value = 152.714640345757
$ws.Cells.Item(1, 5).Value = 43.6933046410295
$ws.Cells.Item(1, 6).Value = 91.3657947868881
$ws.Cells.Item(1, 7).Value = 18.9352595335502
$ws.Cells.Item(1, 8).Value = 132.25148512621
$ws.Cells.Item(1, 9).Value = 187.063584284421
$ws.Cells.Item(1, 10).Value = 6.48669312078817

# Row 2
$ws.Cells.Item(2, 1).Value = 150.866383663782
$ws.Cells.Item(2, 2).Value = 26.0114197740384
$ws.Cells.Item(2, 3).Value = 16.6590415950208
$ws.Cells.Item(2, 4).Value = 72.7183658968277
$ws.Cells.Item(2, 5).Value = 163.46372438756
$ws.Cells.Item(2, 6).Value = 68.4148016704315
$ws.Cells.Item(2, 7).Value = 128.843220290189
$ws.Cells.Item(2, 8).Value = 78.2641375801825
$ws.Cells.Item(2, 9).Value = 174.00358895492
$ws.Cells.Item(2, 10).Value = 30.5255962677885

# Row 3
$ws.Cells.Item(3, 1).Value = 37.0394148105008
$ws.Cells.Item(3, 2).Value = 92.808758417521
$ws.Cells.Item(3, 3).Value = 82.9473020895139
$ws.Cells.Item(3, 4).Value = 13.4338728214772
$ws.Cells.Item(3, 5).Value = 90.14611313592
$ws.Cells.Item(3, 6).Value = 197.089343702928
$ws.Cells.Item(3, 7).Value = 39.1669149692948
$ws.Cells.Item(3, 8).Value = 168.83658048177
$ws.Cells.Item(3, 9).Value = 60.9308183476938
$ws.Cells.Item(3, 10).Value = 65.4367342895999

# Row 4
$ws.Cells.Item(4, 1).Value = 18.0905228564937
$ws.Cells.Item(4, 2).Value = 114.608184860371
$ws.Cells.Item(4, 3).Value = 191.565176281876
$ws.Cells.Item(4, 4).Value = 36.8937738411565
$ws.Cells.Item(4, 5).Value = 146.780665845974
$ws.Cells.Item(4, 6).Value = 38.7618960061864
$ws.Cells.Item(4, 7).Value = 114.890729316925
$ws.Cells.Item(4, 8).Value = 108.486213585588
$ws.Cells.Item(4, 9).Value = 150.807154807638
$ws.Cells.Item(4, 10).Value = 30.9052166672913

# Row 5
$ws.Cells.Item(5, 1).Value = 45.7007184837482
$ws.Cells.Item(5, 2).Value = 163.005481363742
$ws.Cells.Item(5, 3).Value = 169.108958248565
$ws.Cells.Item(5, 4).Value = 34.397506823017
$ws.Cells.Item(5, 5).Value = 174.217484786276
$ws.Cells.Item(5, 6).Value = 151.109676971617
$ws.Cells.Item(5, 7).Value = 31.8791303932104
$ws.Cells.Item(5, 8).Value = 46.2921630806719
$ws.Cells.Item(5, 9).Value = 164.163637517096
$ws.Cells.Item(5, 10).Value = 127.563915647363

# Row 6
$ws.Cells.Item(6, 1).Value = 145.105422542014
$ws.Cells.Item(6, 2).Value = 54.5950544320955
$ws.Cells.Item(6, 3).Value = 193.323494211456
$ws.Cells.Item(6, 4).Value = 142.025130773906
$ws.Cells.Item(6, 5).Value = 7.84210851781168
$ws.Cells.Item(6, 6).Value = 81.8060659253067
$ws.Cells.Item(6, 7).Value = 196.969882583697
$ws.Cells.Item(6, 8).Value = 149.814643873747
$ws.Cells.Item(6, 9).Value = 62.5685272098372
$ws.Cells.Item(6, 10).Value = 46.6039609381016

# Row 7
$ws.Cells.Item(7, 1).Value = 52.1988798175933
$ws.Cells.Item(7, 2).Value = 50.0986790517805
$ws.Cells.Item(7, 3).Value = 71.3206667785163
$ws.Cells.Item(7, 4).Value = 121.626849994821
$ws.Cells.Item(7, 5).Value = 188.396170264294
$ws.Cells.Item(7, 6).Value = 36.258198803411
$ws.Cells.Item(7, 7).Value = 34.4462434921629
$ws.Cells.Item(7, 8).Value = 179.765267753864
$ws.Cells.Item(7, 9).Value = 125.937700050854
$ws.Cells.Item(7, 10).Value = 124.701828381373

# Row 8
$ws.Cells.Item(8, 1).Value = 153.524072353506
$ws.Cells.Item(8, 2).Value = 20.357006704601
$ws.Cells.Item(8, 3).Value = 127.456982772544
$ws.Cells.Item(8, 4).Value = 143.098372287628
$ws.Cells.Item(8, 5).Value = 184.82487778404
$ws.Cells.Item(8, 6).Value = 74.0339334467584
$ws.Cells.Item(8, 7).Value = 123.699800168956
$ws.Cells.Item(8, 8).Value = 48.5497952664969
$ws.Cells.Item(8, 9).Value = 39.2163880352007
$ws.Cells.Item(8, 10).Value = 139.036436164303

# Row 9
$ws.Cells.Item(9, 1).Value = 165.210213309717
$ws.Cells.Item(9, 2).Value = 192.874751888623
$ws.Cells.Item(9, 3).Value = 4.67294296467348
$ws.Cells.Item(9, 4).Value = 133.366902700331
$ws.Cells.Item(9, 5).Value = 120.331311747586
$ws.Cells.Item(9, 6).Value = 163.495468424398
$ws.Cells.Item(9, 7).Value = 121.284690648916
$ws.Cells.Item(9, 8).Value = 49.5400455079694
$ws.Cells.Item(9, 9).Value = 29.0516653233448
$ws.Cells.Item(9, 10).Value = 64.9745999206671

# Row 10
$ws.Cells.Item(10, 1).Value = 41.7920134224892
$ws.Cells.Item(10, 2).Value = 165.076085443178
$ws.Cells.Item(10, 3).Value = 45.917686375751
$ws.Cells.Item(10, 4).Value = 104.203193869536
$ws.Cells.Item(10, 5).Value = 178.706336849698
$ws.Cells.Item(10, 6).Value = 195.602039431968
$ws.Cells.Item(10, 7).Value = 91.684814585226
$ws.Cells.Item(10, 8).Value = 47.4821082537445
$ws.Cells.Item(10, 9).Value = 46.0013365587226
$ws.Cells.Item(10, 10).Value = 137.959285982866

# Row 11
$ws.Cells.Item(11, 1).Value = 116.663433479454
$ws.Cells.Item(11, 2).Value = 52.1138626393461
$ws.Cells.Item(11, 3).Value = 120.354463029818
$ws.Cells.Item(11, 4).Value = 39.4618091357228
$ws.Cells.Item(11, 5).Value = 174.039843293857
$ws.Cells.Item(11, 6).Value = 124.748415837413
$ws.Cells.Item(11, 7).Value = 127.138071659551
$ws.Cells.Item(11, 8).Value = 50.2251219238272
$ws.Cells.Item(11, 9).Value = 157.200252989866
$ws.Cells.Item(11, 10).Value = 133.808175071053

# Row 12
$ws.Cells.Item(12, 1).Value = 158.106265756351
$ws.Cells.Item(12, 2).Value = 148.4200873172
$ws.Cells.Item(12, 3).Value = 110.598255838546
$ws.Cells.Item(12, 4).Value = 123.532091045534
$ws.Cells.Item(12, 5).Value = 81.3937476283841
$ws.Cells.Item(12, 6).Value = 59.3241279289704
$ws.Cells.Item(12, 7).Value = 45.425736087107
$ws.Cells.Item(12, 8).Value = 137.953764078186
$ws.Cells.Item(12, 9).Value = 1.2955382472349
$ws.Cells.Item(12, 10).Value = 24.9007018398962

# Row 13
$ws.Cells.Item(13, 1).Value = 114.973508154495
$ws.Cells.Item(13, 2).Value = 184.906197984194
$ws.Cells.Item(13, 3).Value = 150.713602430519
$ws.Cells.Item(13, 4).Value = 60.9631001301869
$ws.Cells.Item(13, 5).Value = 82.9098149588843
$ws.Cells.Item(13, 6).Value = 188.447986910328
$ws.Cells.Item(13, 7).Value = 174.43932032885
$ws.Cells.Item(13, 8).Value = 23.2537889030081
$ws.Cells.Item(13, 9).Value = 164.39203543793
$ws.Cells.Item(13, 10).Value = 189.222838352072

# Row 14
$ws.Cells.Item(14, 1).Value = 182.349118861532
$ws.Cells.Item(14, 2).Value = 76.2176919152111
$ws.Cells.Item(14, 3).Value = 2.54845870777427
$ws.Cells.Item(14, 4).Value = 101.257102052335
$ws.Cells.Item(14, 5).Value = 22.3730026848488
$ws.Cells.Item(14, 6).Value = 113.096350670371
$ws.Cells.Item(14, 7).Value = 72.520288858805
$ws.Cells.Item(14, 8).Value = 165.211133828951
$ws.Cells.Item(14, 9).Value = 159.327059406474
$ws.Cells.Item(14, 10).Value = 195.582895910173

# Row 15
$ws.Cells.Item(15, 1).Value = 36.3573967648472
$ws.Cells.Item(15, 2).Value = 71.0595687250884
$ws.Cells.Item(15, 3).Value = 92.3397925181034
$ws.Cells.Item(15, 4).Value = 95.2434902522916
$ws.Cells.Item(15, 5).Value = 106.868334164316
$ws.Cells.Item(15, 6).Value = 93.3719261052888
$ws.Cells.Item(15, 7).Value = 54.4778296046321
$ws.Cells.Item(15, 8).Value = 122.385595330217
$ws.Cells.Item(15, 9).Value = 22.8094462411522
$ws.Cells.Item(15, 10).Value = 119.382208920728

# Row 16
$ws.Cells.Item(16, 1).Value = 150.176303344861
$ws.Cells.Item(16, 2).Value = 153.73105050704
$ws.Cells.Item(16, 3).Value = 46.1865700065096
$ws.Cells.Item(16, 4).Value = 21.1006347188263
$ws.Cells.Item(16, 5).Value = 22.9857778283701
$ws.Cells.Item(16, 6).Value = 131.757235495261
$ws.Cells.Item(16, 7).Value = 101.400260208827
$ws.Cells.Item(16, 8).Value = 59.391362899631
$ws.Cells.Item(16, 9).Value = 156.551994176839
$ws.Cells.Item(16, 10).Value = 185.591856383528

# Row 17
$ws.Cells.Item(17, 1).Value = 150.309095508563
$ws.Cells.Item(17, 2).Value = 103.884282756543
$ws.Cells.Item(17, 3).Value = 85.8330864858968
$ws.Cells.Item(17, 4).Value = 167.977414637794
$ws.Cells.Item(17, 5).Value = 151.459056209521
$ws.Cells.Item(17, 6).Value = 81.8885738411399
$ws.Cells.Item(17, 7).Value = 145.871628609426
$ws.Cells.Item(17, 8).Value = 9.34115378621088
$ws.Cells.Item(17, 9).Value = 101.159088360685
$ws.Cells.Item(17, 10).Value = 168.297396958013

# Row 18
$ws.Cells.Item(18, 1).Value = 186.803839070165
$ws.Cells.Item(18, 2).Value = 80.2146022581563
$ws.Cells.Item(18, 3).Value = 178.626704671712
$ws.Cells.Item(18, 4).Value = 5.71264233706176
$ws.Cells.Item(18, 5).Value = 188.543305075049
$ws.Cells.Item(18, 6).Value = 43.913939429407
$ws.Cells.Item(18, 7).Value = 92.5664054660901
$ws.Cells.Item(18, 8).Value = 55.4701121782279
$ws.Cells.Item(18, 9).Value = 154.094765965871
$ws.Cells.Item(18, 10).Value = 189.537888853596

# Row 19
$ws.Cells.Item(19, 1).Value = 133.970157305696
$ws.Cells.Item(19, 2).Value = 52.0537249986332
$ws.Cells.Item(19, 3).Value = 0.444342661855902
$ws.Cells.Item(19, 4).Value = 45.0098265172028
$ws.Cells.Item(19, 5).Value = 39.0465350072116
$ws.Cells.Item(19, 6).Value = 28.6180683544921
$ws.Cells.Item(19, 7).Value = 30.0311219087015
$ws.Cells.Item(19, 8).Value = 181.447823988948
$ws.Cells.Item(19, 9).Value = 78.2713242239651
$ws.Cells.Item(19, 10).Value = 90.6157671895883

# Row 20
$ws.Cells.Item(20, 1).Value = 11.6960904615447
$ws.Cells.Item(20, 2).Value = 13.128925959174
$ws.Cells.Item(20, 3).Value = 8.6591396521121
$ws.Cells.Item(20, 4).Value = 173.735203022945
$ws.Cells.Item(20, 5).Value = 45.2738004016102
$ws.Cells.Item(20, 6).Value = 132.473114008304
$ws.Cells.Item(20, 7).Value = 185.226482239192
$ws.Cells.Item(20, 8).Value = 124.36237788031
$ws.Cells.Item(20, 9).Value = 143.784434042771
$ws.Cells.Item(20, 10).Value = 24.979760323176

Write-Output "Updated A1:J20 with new random values"
